$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the 0.5 value from G7 (Saturday) to H7 (Sunday) for the Read/Study task row.
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 0.5

# Update the active selection on the sheet to H8 (matches the saved view state).
$ws.Range("H8").Select()
